$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3123737
$ws.Range("C3").Value = 1384544
$ws.Range("C4").Value = 23455
$ws.Range("C5").Value = 37224301
$ws.Range("C6").Value = 29399209
$ws.Range("C7").Value = 1178344
$ws.Range("C8").Value = 1057101
$ws.Range("C9").Value = 7802
$ws.Range("C10").Value = 321
$ws.Range("C11").Value = 21969674
$ws.Range("C12").Value = 20661882
$ws.Range("C13").Value = 5602251

$wb.Save()
